$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 7270
$ws.Range("J2").Value = 1051
$ws.Range("J3").Value = 1128
$ws.Range("H4").Value = 1686
$ws.Range("I4").Value = 1751
$ws.Range("J4").Value = 256
$ws.Range("J6").Value = 1537
$ws.Range("H7").Value = 25999
$ws.Range("I7").Value = 26188
$ws.Range("J7").Value = 4055
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 18
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 54
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 14
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 147
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 29
$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 25
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 106
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 22
$ws.Range("J7").Value = 118
$ws.Range("J8").Value = 257
$ws.Range("J11").Value = 49
$ws.Range("J14").Value = 14
$ws.Range("J15").Value = 44
$ws.Range("J16").Value = 9
$ws.Range("J19").Value = 139
$ws.Range("J20").Value = 84
$ws.Range("J24").Value = 17
$ws.Range("J25").Value = 23
$ws.Range("J29").Value = 221
$ws.Range("J31").Value = 29
$ws.Range("I33").Value = 1143
$ws.Range("J33").Value = 165
$ws.Range("J34").Value = 28
$ws.Range("J36").Value = 61
$ws.Range("J42").Value = 174
$ws.Range("J44").Value = 30
$ws.Range("J54").Value = 77
$ws.Range("J60").Value = 24
$ws.Range("H63").Value = 236
$ws.Range("I63").Value = 182
$ws.Range("J63").Value = 19
$ws.Range("J65").Value = 106
$ws.Range("J67").Value = 147
$ws.Range("J76").Value = 63
$ws.Range("J78").Value = 58
$ws.Range("J79").Value = 122
$ws.Range("J85").Value = 166
$ws.Range("J86").Value = 15
$ws.Range("J88").Value = 31
$ws.Range("J92").Value = 13
$ws.Range("J94").Value = 26
$ws.Range("J96").Value = 54
$ws.Range("J97").Value = 22
$ws.Range("J98").Value = 29
$ws.Range("H101").Value = 25999
$ws.Range("I101").Value = 26188
$ws.Range("J101").Value = 4055
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 43
$ws.Range("I4").Value = 49
$ws.Range("J6").Value = 72
$ws.Range("I7").Value = 1143
$ws.Range("J7").Value = 165
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 77
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 62
$ws.Range("J3").Value = 80
$ws.Range("J7").Value = 221
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 31
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 139
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 30
$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 63
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 42
$ws.Range("J3").Value = 58
$ws.Range("J7").Value = 166
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 174
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 14
$ws.Range("J3").Value = 21
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 58
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 6
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 17
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 39
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 122
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 84
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 13
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 61
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 28
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 26
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 23
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 44
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 29
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 49
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 22
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 3
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 13
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 31
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 80
$ws.Range("J3").Value = 86
$ws.Range("J7").Value = 257
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 15
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 6
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 24
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 118
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 22
$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 9
